# Rename the "pairwise_ttests" sheet to "pairwise_tests" and relabel its title,
# and add the new "test__effsize" / "np2" parameter row on the "parameter" sheet.

$wb = $excel.ActiveWorkbook

# 1) Rename the worksheet (tab name)
$pw = $wb.Worksheets.Item("pairwise_ttests")
$pw.Name = "pairwise_tests"

# 2) Update the sheet's own title cell (A1) to match the renamed test
$pw.Range("A1").Value = "Pairwise Tests"

# 3) Append the new parameter row to the "parameter" sheet, matching the
#    bold/bordered label style already used by the other A-column cells.
$paramWs = $wb.Worksheets.Item("parameter")
$paramWs.Range("A6").Copy()
$paramWs.Range("A7").PasteSpecial(-4122)
$paramWs.Range("A7").Value = "test__effsize"
$paramWs.Range("B7").Value = "np2"
